# Update the Naukri test data sheet:
#  - Username cell (C2) text changes from the old gmail address to the
#    new placeholder address (the mailto: hyperlink target itself is left
#    untouched).
#  - Password cell (D2) is re-entered with a leading apostrophe so Excel
#    stores it with an explicit "quote prefix" (text-forced) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = "username@email.com"
$ws.Cells.Item(2, 4).Value = "'password"
